$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the orphaned external reference to the P31 workbook
$wb.BreakLink("/Users/infer/git/alg_AunonAndreaUO277876/src/main/java/algestudiante/p31/P31_UO277876.xlsx", 1)

# Fix header text typos (add missing closing parenthesis)
$ws.Range("B2").Value = "t(O(n2) (ms))"
$ws.Range("C2").Value = "tO(nlogn) (ms))"

# Widen column C slightly to fit the corrected header text
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666

# Move the chart down on the sheet
$co = $ws.ChartObjects(1)
$co.Top = 162.3

# Update the active selection
$ws.Range("G13").Select() | Out-Null
